# The data table on Sheet1 is a quarterly error series keyed by a date
# label in column A. The source was missing the "2020-05-15" observation,
# which belongs right after the existing "2020-02-14" row (row 3) and
# before "2020-08-25" (the old row 3, now pushed down). Insert a new row
# there, which shifts every following row down by one and extends the
# series by one more entry at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 3; rows 3:22 shift down to 4:23.
$ws.Rows("3:3").Insert()

# The inserted row loses the bold/centered/bordered label formatting that
# the rest of column A uses - restore it by copying the format from the
# label cell directly below (A4), then clear the clipboard marquee.
$ws.Range("A4").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the newly inserted observation.
$ws.Range("A3").Value = "2020-05-15 00:00:00_diff"
$ws.Range("B3").Value = -9.678192681000001
$ws.Range("C3").Value = 16.9916
$ws.Range("D3").Value = -4.437406
$ws.Range("E3").Value = 2.102676
$ws.Range("F3").Value = 2.790985
$ws.Range("G3").Value = -0.987134
$ws.Range("H3").Value = -1.382463
